$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.016552413468103
$ws.Range("D2").Value = 1.022040069288516
$ws.Range("E2").Value = 1.044490185947091
$ws.Range("F2").Value = 1.047597120381425
$ws.Range("I2").Value = 1.026437499221819
$ws.Range("J2").Value = 1.021771571777636
$ws.Range("K2").Value = 1.024875623536912
$ws.Range("L2").Value = 1.047261089541029
$ws.Range("M2").Value = 1.050359313613245
$ws.Range("N2").Value = 1.011368696429943
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.017369570045419
$ws.Range("D3").Value = 1.022601070586336
$ws.Range("E3").Value = 1.045709343861217
$ws.Range("F3").Value = 1.048873745858206
$ws.Range("I3").Value = 1.026506656957223
$ws.Range("J3").Value = 1.022224925962693
$ws.Range("K3").Value = 1.025243877923041
$ws.Range("L3").Value = 1.0482904320394
$ws.Range("M3").Value = 1.051446610920002
$ws.Range("N3").Value = 1.011519414690702
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.017898622646614
$ws.Range("D4").Value = 1.022964095842561
$ws.Range("E4").Value = 1.04649914248666
$ws.Range("F4").Value = 1.04970078658907
$ws.Range("I4").Value = 1.026550015102095
$ws.Range("J4").Value = 1.022517938126707
$ws.Range("K4").Value = 1.025481480258419
$ws.Range("L4").Value = 1.048956842703994
$ws.Range("M4").Value = 1.052150590114621
$ws.Range("N4").Value = 1.011616797556645
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018121106419786
$ws.Range("D5").Value = 1.023116714876444
$ws.Range("E5").Value = 1.046831395374309
$ws.Range("F5").Value = 1.050048709619265
$ws.Range("I5").Value = 1.026567909172561
$ws.Range("J5").Value = 1.022641038573422
$ws.Range("K5").Value = 1.025581203770746
$ws.Range("L5").Value = 1.049237087947286
$ws.Range("M5").Value = 1.052446645339062
$ws.Range("N5").Value = 1.011657703008644
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.018158466497756
$ws.Range("D6").Value = 1.02314234042725
$ws.Range("E6").Value = 1.046887195112073
$ws.Range("F6").Value = 1.050107141237455
$ws.Range("I6").Value = 1.026570894081555
$ws.Range("J6").Value = 1.022661702853899
$ws.Range("K6").Value = 1.025597938110434
$ws.Range("L6").Value = 1.04928414743141
$ws.Range("M6").Value = 1.052496360334465
$ws.Range("N6").Value = 1.011664569192846
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.017901595211896
$ws.Range("D7").Value = 1.022966135134882
$ws.Range("E7").Value = 1.046503581194424
$ws.Range("F7").Value = 1.049705434628789
$ws.Range("I7").Value = 1.026550255515426
$ws.Range("J7").Value = 1.022519583323683
$ws.Range("K7").Value = 1.025482813416578
$ws.Range("L7").Value = 1.048960587015085
$ws.Range("M7").Value = 1.052154545616016
$ws.Range("N7").Value = 1.011617344272504
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.016828513007022
$ws.Range("D8").Value = 1.022229657253031
$ws.Range("E8").Value = 1.044902015679088
$ws.Range("F8").Value = 1.048028359780086
$ws.Range("I8").Value = 1.026461159193722
$ws.Range("J8").Value = 1.021924854470386
$ws.Range("K8").Value = 1.025000217608229
$ws.Range("L8").Value = 1.047608887094199
$ws.Range("M8").Value = 1.050726683157676
$ws.Range("N8").Value = 1.011419661562277
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.014939937594953
$ws.Range("D9").Value = 1.020932109775352
$ws.Range("E9").Value = 1.042086876261905
$ws.Range("F9").Value = 1.045080598977255
$ws.Range("I9").Value = 1.026293527180697
$ws.Range("J9").Value = 1.020874313726644
$ws.Range("K9").Value = 1.024144637346443
$ws.Range("L9").Value = 1.045229724885357
$ws.Range("M9").Value = 1.048213836459047
$ws.Range("N9").Value = 1.011070246775643
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.013682525836241
$ws.Range("D10").Value = 1.020067318076599
$ws.Range("E10").Value = 1.040214784779581
$ws.Range("F10").Value = 1.043120387512676
$ws.Range("I10").Value = 1.026174656491114
$ws.Range("J10").Value = 1.020172288337255
$ws.Range("K10").Value = 1.023570825079615
$ws.Range("L10").Value = 1.043645391706557
$ws.Range("M10").Value = 1.046540737204421
$ws.Range("N10").Value = 1.010836601628473
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.013138455780668
$ws.Range("D11").Value = 1.019692929849544
$ws.Range("E11").Value = 1.039405242991294
$ws.Range("F11").Value = 1.042272758839445
$ws.Range("I11").Value = 1.026121503220424
$ws.Range("J11").Value = 1.01986791929926
$ws.Range("K11").Value = 1.023321558184087
$ws.Range("L11").Value = 1.042959768479899
$ws.Range("M11").Value = 1.045816762947645
$ws.Range("N11").Value = 1.010735268183661
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01293642478953
$ws.Range("D12").Value = 1.019553877790871
$ws.Range("E12").Value = 1.03910470513312
$ws.Range("F12").Value = 1.04195808429557
$ws.Range("I12").Value = 1.026101507635593
$ws.Range("J12").Value = 1.019754805706268
$ws.Range("K12").Value = 1.02322884980166
$ws.Range("L12").Value = 1.04270515712393
$ws.Range("M12").Value = 1.045547919201015
$ws.Range("N12").Value = 1.010697604169045
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.012979758351148
$ws.Range("D13").Value = 1.019583704343786
$ws.Range("E13").Value = 1.039169164211611
$ws.Range("F13").Value = 1.042025575266681
$ws.Range("I13").Value = 1.02610580815568
$ws.Range("J13").Value = 1.019779071539304
$ws.Range("K13").Value = 1.023248741458785
$ws.Range("L13").Value = 1.042759769441063
$ws.Range("M13").Value = 1.04560558382971
$ws.Range("N13").Value = 1.01070568432093
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.013121754588029
$ws.Range("D14").Value = 1.019681435501197
$ws.Range("E14").Value = 1.039380397140088
$ws.Range("F14").Value = 1.042246744231528
$ws.Range("I14").Value = 1.026119855518713
$ws.Range("J14").Value = 1.019858570460158
$ws.Range("K14").Value = 1.023313897313511
$ws.Range("L14").Value = 1.042938721014855
$ws.Range("M14").Value = 1.04579453876535
$ws.Range("N14").Value = 1.010732155357593
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.013209251253274
$ws.Range("D15").Value = 1.019741652581642
$ws.Range("E15").Value = 1.039510566151221
$ws.Range("F15").Value = 1.04238303655165
$ws.Range("I15").Value = 1.026128477173702
$ws.Range("J15").Value = 1.019907544790269
$ws.Range("K15").Value = 1.023354026172947
$ws.Range("L15").Value = 1.04304898685892
$ws.Range("M15").Value = 1.045910969721152
$ws.Range("N15").Value = 1.010748461831354
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.013718642445773
$ws.Range("D16").Value = 1.020092166644641
$ws.Range("E16").Value = 1.040268534185624
$ws.Range("F16").Value = 1.043176666088556
$ws.Range("I16").Value = 1.026178148717048
$ws.Range("J16").Value = 1.020192480215946
$ws.Range("K16").Value = 1.023587351281714
$ws.Range("L16").Value = 1.043690902679204
$ws.Range("M16").Value = 1.046588795211859
$ws.Range("N16").Value = 1.010843323376126
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.014038277309733
$ws.Range("D17").Value = 1.020312055296754
$ws.Range("E17").Value = 1.040744277474848
$ws.Range("F17").Value = 1.043674797553845
$ws.Range("I17").Value = 1.026208856558618
$ws.Range("J17").Value = 1.020371109541748
$ws.Range("K17").Value = 1.023733495835042
$ws.Range("L17").Value = 1.044093666845652
$ws.Range("M17").Value = 1.047014107450042
$ws.Range("N17").Value = 1.010902783948581
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.014224753198719
$ws.Range("D18").Value = 1.020440319498954
$ws.Range("E18").Value = 1.041021875365607
$ws.Range("F18").Value = 1.043965460808001
$ws.Range("I18").Value = 1.026226605649377
$ws.Range("J18").Value = 1.020475263570546
$ws.Range("K18").Value = 1.023818661994334
$ws.Range("L18").Value = 1.044328631366663
$ws.Range("M18").Value = 1.047262232101487
$ws.Range("N18").Value = 1.010937450493094
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.014288343116656
$ws.Range("D19").Value = 1.020484055395914
$ws.Range("E19").Value = 1.04111654689786
$ws.Range("F19").Value = 1.044064588430744
$ws.Range("I19").Value = 1.026232630090095
$ws.Range("J19").Value = 1.020510771028103
$ws.Range("K19").Value = 1.023847688280541
$ws.Range("L19").Value = 1.04440875484085
$ws.Range("M19").Value = 1.047346844232192
$ws.Range("N19").Value = 1.010949268199245
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.014003979546102
$ws.Range("D20").Value = 1.020288462603492
$ws.Range("E20").Value = 1.040693223882508
$ws.Range("F20").Value = 1.043621341226991
$ws.Range("I20").Value = 1.026205578678732
$ws.Range("J20").Value = 1.020351948171696
$ws.Range("K20").Value = 1.023717823915131
$ws.Range("L20").Value = 1.044050450035621
$ws.Range("M20").Value = 1.046968470594326
$ws.Range("N20").Value = 1.010896406024045
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.013079938546557
$ws.Range("D21").Value = 1.01965265575664
$ws.Range("E21").Value = 1.039318189844746
$ws.Range("F21").Value = 1.042181610731725
$ws.Range("I21").Value = 1.026115725875185
$ws.Range("J21").Value = 1.019835161586792
$ws.Range("K21").Value = 1.023294713830716
$ws.Range("L21").Value = 1.042886022598189
$ws.Range("M21").Value = 1.045738894251498
$ws.Range("N21").Value = 1.010724360964309
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.012499309925464
$ws.Range("D22").Value = 1.019252971247789
$ws.Range("E22").Value = 1.03845458811309
$ws.Range("F22").Value = 1.041277392406397
$ws.Range("I22").Value = 1.026057773578082
$ws.Range("J22").Value = 1.019509905771843
$ws.Range("K22").Value = 1.023027996485932
$ws.Range("L22").Value = 1.042154245416275
$ws.Range("M22").Value = 1.044966229980812
$ws.Range("N22").Value = 1.010616049191351
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.012807078341612
$ws.Range("D23").Value = 1.019464844184859
$ws.Range("E23").Value = 1.038912311330606
$ws.Range("F23").Value = 1.041756641495918
$ws.Range("I23").Value = 1.02608863323807
$ws.Range("J23").Value = 1.019682361188478
$ws.Range("K23").Value = 1.023169453588059
$ws.Range("L23").Value = 1.042542142011998
$ws.Range("M23").Value = 1.045375794496973
$ws.Range("N23").Value = 1.010673480498789
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.014019477110968
$ws.Range("D24").Value = 1.020299123104983
$ws.Range("E24").Value = 1.040716292478891
$ws.Range("F24").Value = 1.043645495494909
$ws.Range("I24").Value = 1.02620706031301
$ws.Range("J24").Value = 1.020360606485767
$ws.Range("K24").Value = 1.023724905620665
$ws.Range("L24").Value = 1.044069977731019
$ws.Range("M24").Value = 1.046989091778802
$ws.Range("N24").Value = 1.01089928798246
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.015427895951599
$ws.Range("D25").Value = 1.021267521291742
$ws.Range("E25").Value = 1.042813830546053
$ws.Range("F25").Value = 1.045841787298414
$ws.Range("I25").Value = 1.026338120730847
$ws.Range("J25").Value = 1.021146200886583
$ws.Range("K25").Value = 1.024366433527624
$ws.Range("L25").Value = 1.04584447927806
$ws.Range("M25").Value = 1.048863088942891
$ws.Range("N25").Value = 1.011160703999875